# Updates the cryptos price/volume sheet with refreshed market data.
# Price (column D) and Volume(1h) (column E) values are stored as plain
# text in this sheet (e.g. "27.863.80", "  -0.93%  "), so for any D-column
# value that looks like a genuine number we force NumberFormat to "@"
# (Text) before assigning it, otherwise Excel's COM layer would silently
# coerce strings like "1.00" or "0.999" into numeric values and lose the
# original text representation/formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.863.80"
$ws.Range("D3").Value = "1.624.44"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.97"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.41"
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.257"
$ws.Range("E9").Value = "  -1.82%  "
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0879"
$ws.Range("E11").Value = "  -0.11%  "
$ws.Range("D12").Value = "1.856.00"
$ws.Range("E12").Value = "  -0.88%  "
$ws.Range("D13").Value = "1.621.81"
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("E15").Value = "  -2.42%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.31"
$ws.Range("E16").Value = "  -0.90%  "
$ws.Range("D17").Value = "27.854.27"
$ws.Range("E17").Value = "  -0.41%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "229.70"
$ws.Range("E18").Value = "  -1.43%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.63"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "0.0₃0720"
$ws.Range("E20").Value = "  -0.65%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("E22").Value = "  -1.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.13"
$ws.Range("E23").Value = "  -5.18%  "
$ws.Range("E24").Value = "  -2.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.99"
$ws.Range("E25").Value = "  +2.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.89"
$ws.Range("E26").Value = "  -1.53%  "
$ws.Range("E27").Value = "  -0.28%  "
$ws.Range("E28").Value = "  -1.20%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("E31").Value = "  -0.57%  "
$ws.Range("E32").Value = "  +2.11%  "
$ws.Range("E33").Value = "  -1.71%  "
$ws.Range("D34").Value = "1.391.80"
$ws.Range("E34").Value = "  -1.07%  "
$ws.Range("E35").Value = "  -0.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.01"
$ws.Range("E36").Value = "  +11.86%  "
$ws.Range("E37").Value = "  -1.05%  "
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.859"
$ws.Range("E40").Value = "  -2.61%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.16%  "
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.83"
$ws.Range("E44").Value = "  -3.02%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.64"
$ws.Range("E45").Value = "  -1.05%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.18"
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("D47").Value = "1.765.93"
$ws.Range("E47").Value = "  -0.85%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.84"
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0105"
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.102"
$ws.Range("E50").Value = "  +1.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0503"
$ws.Range("E51").Value = "  -0.50%  "
